# Part1-EventBasedGameLoop.pptx edit:
# Slide 1 ("Today's Attendance password" textbox) — the blank password
# line ("_________") is filled in with the answer "timers".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 4 on slide 1 is "TextBox 2" (the attendance-password callout),
# whose 2nd paragraph currently just holds the blank "_________".
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(2)
$para.Text = "timers"
